# Update the weekly schedule workbook:
#  - refresh several of the "assigned person" cells in column B (rows 1-6)
#    on each day's sheet (some people reassigned, some cells cleared out)
#  - add new "Start time:" / "End time:" rows (7 and 8) to every sheet,
#    capturing the shift hours for that day (still a work in progress per
#    the commit message, so some sheets only get partial data)

$wb = $excel.ActiveWorkbook

function Set-DaySheet {
    param(
        [string]$SheetName,
        [hashtable]$Assignments,   # row number (1-6) -> name, or $null to clear
        $StartTime,                # value for B7, or $null to leave blank
        $EndTime                   # value for B8, or $null to leave blank
    )

    $ws = $wb.Worksheets.Item($SheetName)
    # NOTE: this runtime's PowerShell only reliably binds *positional*
    # arguments to user-defined functions, so callers below pass
    # arguments positionally (no -Name style parameters).

    foreach ($row in 1..6) {
        $value = $Assignments[$row]
        $cell = $ws.Range("B$row")
        if ($null -eq $value) {
            $cell.ClearContents()
        } else {
            $cell.Value = $value
        }
    }

    $ws.Range("A7").Value = "Start time: "
    if ($null -ne $StartTime) {
        $ws.Range("B7").Value = $StartTime
    }

    $ws.Range("A8").Value = "End time: "
    if ($null -ne $EndTime) {
        $ws.Range("B8").Value = $EndTime
    }
}

# Monday
Set-DaySheet "Monday" @{
    1 = "fred york"
    2 = "Emmie"
    3 = "flora"
    4 = "Jordan"
    5 = "Jordan"
    6 = "Ruby"
} 700 1500

# Tuesday
Set-DaySheet "Tuesday" @{
    1 = "Emmie"
    2 = "Ruby"
    3 = "flora"
    4 = "Ruby"
    5 = "Ruby"
    6 = "Emmie"
} 1500 2200

# Wednesday
Set-DaySheet "Wednesday" @{
    1 = $null
    2 = "Jordan"
    3 = $null
    4 = "Ruby"
    5 = "Jordan"
    6 = $null
} 600 1500

# Thursday
Set-DaySheet "Thursday" @{
    1 = "Emmie"
    2 = "Jordan"
    3 = "Emmie"
    4 = "Jordan"
    5 = "Jordan"
    6 = "Emmie"
} 1400 2300

# Friday
Set-DaySheet "Friday" @{
    1 = "Emmie"
    2 = "Jordan"
    3 = "Emmie"
    4 = "Jordan"
    5 = "Jordan"
    6 = "Emmie"
} 1400 2300

# Saturday
Set-DaySheet "Saturday" @{
    1 = $null
    2 = $null
    3 = $null
    4 = $null
    5 = $null
    6 = $null
} 600 1400

# Sunday
Set-DaySheet "Sunday" @{
    1 = $null
    2 = $null
    3 = $null
    4 = $null
    5 = $null
    6 = $null
} $null $null
